# Marksheet fix: handle the float-grading-input case without breaking the
# rest of the sheet (was showing "Absent" with all-zero tallies even though
# the student actually answered the quiz). Also trims the two extra
# (empty/duplicate) answer-key blocks that used to live in columns D:E and
# G:H, now that only one "Student Ans / Correct Ans" block is needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the 2nd/3rd answer-key blocks -------------------------------
# Block 3 (G:H) is removed entirely for rows 15-21 (the only rows it used).
$ws.Range("G15:H21").Clear()
# Block 2 (D:E) keeps its header + first three data rows (16-18); the rest
# (19-40) is removed.
$ws.Range("D19:E40").Clear()

# --- 2. Give the row-label cells (No. / Marking / Total) the same header
#        style used elsewhere in that block (copy format only). -----------
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)

# --- 3. Fill in the real tally numbers now that the student isn't marked
#        "Absent" any more. ------------------------------------------------
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 19
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
# C11 used to hold "-1" as text (the bug this commit fixes); store it as a
# real number instead.
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 32
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "31/112"

# --- 4. Populate column A with the student's actual picked answer for the
#        questions that were attempted, colour-coded green/red to match the
#        existing correct/incorrect styles (copied from the existing
#        green B10:B12 / red C10:C12 cells so the shared style indices are
#        reused rather than duplicated). ------------------------------------
function Set-AnswerCell([string]$addr, [string]$text, [bool]$correct) {
    if ($correct) {
        $ws.Range("B10").Copy()
    } else {
        $ws.Range("C10").Copy()
    }
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $text
}

Set-AnswerCell "A16" "Option A" $true
Set-AnswerCell "D16" "Option A" $true
Set-AnswerCell "A18" "Option B" $true
Set-AnswerCell "D18" "Option D" $true
Set-AnswerCell "A19" "Option C" $true
Set-AnswerCell "A27" "Option A" $true
Set-AnswerCell "A31" "Option B" $false
Set-AnswerCell "A37" "Option A" $true
Set-AnswerCell "A38" "Option A" $true

Write-Output "done"
